$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1648.4642
$ws.Range("I15").Value = 1648.4642
$ws.Range("K15").Value = 4945.392599999999
$ws.Range("M15").Value = -4776.392599999999

$ws.Range("H40").Value = 7087.4165
$ws.Range("I40").Value = 5079
$ws.Range("K40").Value = 5079
$ws.Range("M40").Value = -4904

$ws.Range("H62").Value = 4869
$ws.Range("I62").Value = 2440.4
$ws.Range("K62").Value = 2440.4
$ws.Range("M62").Value = -1816.4

$ws.Range("H65").Value = 4869
$ws.Range("I65").Value = 2440.4
$ws.Range("K65").Value = 12202
$ws.Range("M65").Value = -9082

$ws.Range("H70").Value = 4581.2856
$ws.Range("I70").Value = 2397.5
$ws.Range("J70").Value = 6219.125
$ws.Range("K70").Value = 7192.5
$ws.Range("L70").Value = 18657.375
$ws.Range("M70").Value = -6922.5
$ws.Range("N70").Value = -19197.375

$ws.Range("H73").Value = 4581.2856
$ws.Range("I73").Value = 2397.5
$ws.Range("J73").Value = 6219.125
$ws.Range("K73").Value = 7192.5
$ws.Range("L73").Value = 18657.375
$ws.Range("M73").Value = -6256.5
$ws.Range("N73").Value = -20529.375

$ws.Range("H100").Value = 1879.8
$ws.Range("I100").Value = 1974.75
$ws.Range("K100").Value = 1974.75
$ws.Range("M100").Value = -1433.75

$ws.Range("H112").Value = 2333.7144
$ws.Range("J112").Value = 1767.2
$ws.Range("L112").Value = 5301.6
$ws.Range("N112").Value = -7517.6

$ws.Range("H135").Value = 694.0714
$ws.Range("I135").Value = 694.0714
$ws.Range("K135").Value = 6246.6426
$ws.Range("M135").Value = -3711.6426

$ws.Range("H137").Value = 1977.0555
$ws.Range("I137").Value = 878.9
$ws.Range("J137").Value = 3349.75
$ws.Range("K137").Value = 2636.7
$ws.Range("L137").Value = 10049.25
$ws.Range("M137").Value = -86.69999999999982
$ws.Range("N137").Value = -15149.25

$ws.Range("H138").Value = 2314.7856
$ws.Range("J138").Value = 6072.3335
$ws.Range("L138").Value = 18217.0005
$ws.Range("N138").Value = -28497.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2988.1177
$ws.Range("I45").Value = 1961.3334
$ws.Range("K45").Value = 1961.3334
$ws.Range("M45").Value = -1584.3334

$ws.Range("H74").Value = 4007.6667
$ws.Range("I74").Value = 3407.6667
$ws.Range("K74").Value = 3407.6667
$ws.Range("M74").Value = -2533.6667

$ws.Range("H77").Value = 4007.6667
$ws.Range("I77").Value = 3407.6667
$ws.Range("K77").Value = 17038.3335
$ws.Range("M77").Value = -12670.3335

$ws.Range("H110").Value = 100001600
$ws.Range("I110").Value = 142858450
$ws.Range("J110").Value = 2304.3333
$ws.Range("K110").Value = 142858450
$ws.Range("L110").Value = 2304.3333
$ws.Range("M110").Value = -142856405
$ws.Range("N110").Value = -6394.3333

$ws.Range("H122").Value = 2475
$ws.Range("I122").Value = 2475
$ws.Range("K122").Value = 7425
$ws.Range("M122").Value = -4975

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1398.6
$ws.Range("I94").Value = 1398.6
$ws.Range("K94").Value = 1398.6
$ws.Range("M94").Value = -947.5999999999999

$ws.Range("H107").Value = 50007124
$ws.Range("I107").Value = 100004260
$ws.Range("J107").Value = 9998
$ws.Range("K107").Value = 100004260
$ws.Range("L107").Value = 9998
$ws.Range("M107").Value = -100002340
$ws.Range("N107").Value = -13838

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5524.636
$ws.Range("I31").Value = 1615.8889
$ws.Range("K31").Value = 1615.8889
$ws.Range("M31").Value = -1320.8889

$ws.Range("H34").Value = 5524.636
$ws.Range("I34").Value = 1615.8889
$ws.Range("K34").Value = 1615.8889
$ws.Range("M34").Value = -1413.8889

$ws.Range("H58").Value = 3226.4119
$ws.Range("J58").Value = 6321.1665
$ws.Range("L58").Value = 6321.1665
$ws.Range("N58").Value = -6727.1665

$ws.Range("H80").Value = 40000
$ws.Range("J80").Value = 40000
$ws.Range("L80").Value = 40000
$ws.Range("N80").Value = -42246

$ws.Range("H83").Value = 40000
$ws.Range("J83").Value = 40000
$ws.Range("L83").Value = 120000
$ws.Range("N83").Value = -131232

$ws.Range("H107").Value = 1264.75
$ws.Range("I107").Value = 385.22223
$ws.Range("J107").Value = 2395.5715
$ws.Range("K107").Value = 385.22223
$ws.Range("L107").Value = 2395.5715
$ws.Range("M107").Value = 1534.77777
$ws.Range("N107").Value = -6235.5715

$ws.Range("H132").Value = 1672.2
$ws.Range("I132").Value = 1672.2
$ws.Range("K132").Value = 5016.6
$ws.Range("M132").Value = -2486.6

$ws.Range("H136").Value = 3226.4119
$ws.Range("J136").Value = 6321.1665
$ws.Range("L136").Value = 18963.4995
$ws.Range("N136").Value = -24063.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 296884.4
$ws.Range("I128").Value = 296884.4
$ws.Range("K128").Value = 890653.2000000001
$ws.Range("M128").Value = -885673.2000000001

$ws.Range("H129").Value = 2463.3
$ws.Range("I129").Value = 1200.5
$ws.Range("J129").Value = 2779
$ws.Range("K129").Value = 3601.5
$ws.Range("L129").Value = 8337
$ws.Range("M129").Value = 1398.5
$ws.Range("N129").Value = -18337

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1182.8334
$ws.Range("J80").Value = 1132.6666
$ws.Range("L80").Value = 1132.6666
$ws.Range("N80").Value = -3128.6666

$ws.Range("H83").Value = 1182.8334
$ws.Range("J83").Value = 1132.6666
$ws.Range("L83").Value = 5663.333000000001
$ws.Range("N83").Value = -15647.333

$ws.Range("H122").Value = 4399.2
$ws.Range("I122").Value = 3999
$ws.Range("K122").Value = 11997
$ws.Range("M122").Value = -9547

$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820

$ws.Range("H125").Value = 25000
$ws.Range("J125").Value = 25000
$ws.Range("L125").Value = 25000
$ws.Range("N125").Value = -29920

$ws.Range("H130").Value = 50000
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040

$ws.Range("H132").Value = 84742.836
$ws.Range("I132").Value = 92365.45
$ws.Range("K132").Value = 277096.35
$ws.Range("M132").Value = -274566.35

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 83337590
$ws.Range("I61").Value = 125002380
$ws.Range("K61").Value = 125002380
$ws.Range("M61").Value = -125002178

$ws.Range("H80").Value = 75750
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 75750
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H113").Value = 83337590
$ws.Range("I113").Value = 125002380
$ws.Range("K113").Value = 125002380
$ws.Range("M113").Value = -125000210

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 25313
$ws.Range("J74").Value = 25313
$ws.Range("L74").Value = 25313
$ws.Range("N74").Value = -27185

$ws.Range("H77").Value = 25313
$ws.Range("J77").Value = 25313
$ws.Range("L77").Value = 75939
$ws.Range("N77").Value = -85299

$ws.Range("H99").Value = 44000
$ws.Range("J99").Value = 44000
$ws.Range("L99").Value = 44000
$ws.Range("N99").Value = -49990

$ws.Range("H107").Value = 30303646
$ws.Range("I107").Value = 41666990
$ws.Range("K107").Value = 125000970
$ws.Range("M107").Value = -124999050

$ws.Range("H122").Value = 1578.9524
$ws.Range("I122").Value = 1466.25
$ws.Range("K122").Value = 4398.75
$ws.Range("M122").Value = -1948.75

$ws.Range("H132").Value = 1553.5454
$ws.Range("I132").Value = 1553.5454
$ws.Range("K132").Value = 4660.6362
$ws.Range("M132").Value = -2130.6362
